# SpartaFarming ItemInfo.xlsx update
# [feat] UI 추가 및 ItemObject Prefab 추가
#
# - Adds a new "Beet" seed row (inserted as row 5)
# - Repurposes former "CarrotSeed" row (row 4) into a "Wheat" seed row
# - Adds a new "price" / "Price" (int) column (column M)
# - Updates maxStack (L) values for the seed/resource/food rows to 999
# - Adds price values for every item row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row for "Beet" right after the existing seed row (row 4),
#    then remove one of the two duplicate blank spacer rows further down
#    (rows 10/11 in the original sheet were identical blank rows) so the
#    overall row count of the sheet stays the same as before.
# ---------------------------------------------------------------------
$ws.Rows("5").Insert()
$ws.Rows("11").Delete()

# ---------------------------------------------------------------------
# 2. Row 5 (newly inserted) -> "Beet"  (filled in first so its strings
#    are registered before Wheat's, matching shared-string ordering)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Beet"
$ws.Range("C5").Value = "비트"
$ws.Range("D5").Value = "Seed"
$ws.Range("E5").Value = "Sprite/Seed/BeetSeed"
$ws.Range("F5").Value = "Prefab/Seed/BeetSeed"
$ws.Range("G5").Value = -1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 999

# ---------------------------------------------------------------------
# 3. Row 4 used to describe "CarrotSeed" -> turn it into "Wheat"
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "Wheat"
$ws.Range("C4").Value = "밀"
$ws.Range("D4").Value = "Seed"
$ws.Range("E4").Value = "Sprite/Seed/WheatSeed"
$ws.Range("F4").Value = "Prefab/Seed/WheatSeed"
$ws.Range("L4").Value = 999

# ---------------------------------------------------------------------
# 4. maxStack (L) updates for the rows that shifted down (Log, Bread)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 999
$ws.Range("L7").Value = 999

# ---------------------------------------------------------------------
# 5. New "price" column (M) - header + data
# ---------------------------------------------------------------------
# Copy the styling from column L's header cells onto column M's so the
# new column blends in with the rest of the table header.
$ws.Range("L1:L3").Copy()
$ws.Range("M1:M3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M1").Value = "price"
$ws.Range("M2").Value = "int"
$ws.Range("M3").Value = "Price"

$ws.Range("M4").Value = 20
$ws.Range("M5").Value = 60
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 100
$ws.Range("M8").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("M10").Value = 0

# ---------------------------------------------------------------------
# 6. Misc bookkeeping to mirror the authored workbook
# ---------------------------------------------------------------------
$ws.Range("N11").Select()
